$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style (incl. date number format) of the last existing date cell (A366)
# so the new date cells re-use the same style index instead of creating new ones.
$ws.Range("A366").Copy()
$ws.Range("A367:A374").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$data = @(
    @(367, 44441, 7,  43, 130.3701907043022),
    @(368, 44442, 2,  33, 100.0515417033017),
    @(369, 44443, 3,  35, 106.1152715035018),
    @(370, 44444, 6,  28, 84.89221720280145),
    @(371, 44445, 10, 28, 84.89221720280145),
    @(372, 44446, 3,  31, 93.9878119031016),
    @(373, 44447, 1,  32, 97.01967680320165),
    @(374, 44448, 2,  27, 81.8603523027014)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
